# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp note in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 17:14"

# --- Swap Groenlandia / Islas Malvinas order: row 209 becomes Groenlandia,
#     row 210 becomes Islas Malvinas ---
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"

# --- Update per-country statistics (B:Casos totales, C:Nuevos casos,
#     D:Casos activos, E:Recuperados, F:Casos criticos, G:Muertes hoy,
#     H:Muertes) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 3430726
$ws.Cells.Item(4, 3).Value = 16731
$ws.Cells.Item(4, 4).Value = 1518254
$ws.Cells.Item(4, 5).Value = 1774633
$ws.Cells.Item(4, 7).Value = 57
$ws.Cells.Item(4, 8).Value = 137839

# Row 6: India
$ws.Cells.Item(6, 2).Value = 898680
$ws.Cells.Item(6, 3).Value = 19214
$ws.Cells.Item(6, 4).Value = 566664
$ws.Cells.Item(6, 5).Value = 308447
$ws.Cells.Item(6, 7).Value = 382
$ws.Cells.Item(6, 8).Value = 23569

# Row 12: Reino Unido
$ws.Cells.Item(12, 2).Value = 290133
$ws.Cells.Item(12, 3).Value = 530
$ws.Cells.Item(12, 7).Value = 11
$ws.Cells.Item(12, 8).Value = 44830

# Row 19: Alemania
$ws.Cells.Item(19, 2).Value = 200180
$ws.Cells.Item(19, 3).Value = 230
$ws.Cells.Item(19, 5).Value = 5943
$ws.Cells.Item(19, 7).Value = 3
$ws.Cells.Item(19, 8).Value = 9137

# Row 45: Republica Dominicana
$ws.Cells.Item(45, 2).Value = 45506
$ws.Cells.Item(45, 3).Value = 974
$ws.Cells.Item(45, 4).Value = 22441
$ws.Cells.Item(45, 5).Value = 22162
$ws.Cells.Item(45, 7).Value = 6
$ws.Cells.Item(45, 8).Value = 903

# Row 46: Israel
$ws.Cells.Item(46, 2).Value = 39979
$ws.Cells.Item(46, 3).Value = 1309
$ws.Cells.Item(46, 4).Value = 19282
$ws.Cells.Item(46, 5).Value = 20333

# Row 89: Tayikistan
$ws.Cells.Item(89, 2).Value = 6596
$ws.Cells.Item(89, 3).Value = 44
$ws.Cells.Item(89, 4).Value = 5278
$ws.Cells.Item(89, 5).Value = 1263

# Row 90: Estado de Palestina
$ws.Cells.Item(90, 2).Value = 6566
$ws.Cells.Item(90, 3).Value = 336
$ws.Cells.Item(90, 5).Value = 5444

# Row 96: Republica de Yibuti
$ws.Cells.Item(96, 2).Value = 4977
$ws.Cells.Item(96, 3).Value = 5
$ws.Cells.Item(96, 4).Value = 4729
$ws.Cells.Item(96, 5).Value = 192

# Row 110: Sri Lanka
$ws.Cells.Item(110, 2).Value = 2642
$ws.Cells.Item(110, 3).Value = 25
$ws.Cells.Item(110, 5).Value = 650

# Row 111: Cuba
$ws.Cells.Item(111, 2).Value = 2428
$ws.Cells.Item(111, 3).Value = 2
$ws.Cells.Item(111, 4).Value = 2268
$ws.Cells.Item(111, 5).Value = 73

# Row 142: Republica de Chipre
$ws.Cells.Item(142, 2).Value = 1022
$ws.Cells.Item(142, 3).Value = 1
$ws.Cells.Item(142, 5).Value = 164

# Row 156: Reunion
$ws.Cells.Item(156, 2).Value = 596
$ws.Cells.Item(156, 3).Value = 3
$ws.Cells.Item(156, 5).Value = 121

# Row 170: Lesoto
$ws.Cells.Item(170, 4).Value = 33
$ws.Cells.Item(170, 5).Value = 209
$ws.Cells.Item(170, 7).Value = 1
$ws.Cells.Item(170, 8).Value = 3
